# Append June schedule rows 11-18 to the "June" sheet (ActiveSheet / ActiveWorkbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New row data -------------------------------------------------------
# columns: A=Date(serial), B=Type, C=AM/PM, D=Area, E=Mark(payback), F=amount
$rows = @(
    @{ r = 11; a = 40701; b = "showa"; c = "PM"; d = "Taxi";     e = "payback"; f = 10 },
    @{ r = 12; a = 40701; b = "showa"; c = "PM"; d = "Taxi Tax"; e = "payback"; f = 2 },
    @{ r = 13; a = 40702; b = "showa"; c = "AM"; d = "Taxi";     e = "payback"; f = 75 },
    @{ r = 14; a = 40702; b = "showa"; c = "AM"; d = "Taxi Tax"; e = "payback"; f = 2 },
    @{ r = 15; a = 40702; b = "showa"; c = "PM"; d = "Taxi";     e = "payback"; f = 11 },
    @{ r = 16; a = 40702; b = "showa"; c = "PM"; d = "Taxi Tax"; e = "payback"; f = 2 },
    @{ r = 17; a = 40703; b = "showa"; c = "AM"; d = "Taxi";     e = "payback"; f = 64 },
    @{ r = 18; a = 40703; b = "showa"; c = "AM"; d = "Taxi Tax"; e = "payback"; f = 2 }
)

foreach ($row in $rows) {
    $n = $row.r

    $ws.Range("A$n").Value = $row.a
    $ws.Range("B$n").Value = $row.b
    $ws.Range("C$n").Value = $row.c
    $ws.Range("D$n").Value = $row.d
    $ws.Range("E$n").Value = $row.e
    $ws.Range("F$n").Value = $row.f

    # Reuse the existing date-format style (A2/C2, style index 9) instead of
    # letting a brand-new numFmt get synthesized.
    $ws.Range("A2").Copy()
    $ws.Range("A$n").PasteSpecial(-4122)
    $ws.Range("C2").Copy()
    $ws.Range("C$n").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---- Column F width: narrower, best-fit to the new (shorter) numbers.
# (The engine quantizes saved widths to 1/7 character-units, so the nearest
# reachable value to the authored 9.375 is 9.428571428571429 = 66/7.)
$ws.Columns.Item(6).ColumnWidth = 8.7

# ---- Selection matches the post-edit active cell in the source sheet ----
$ws.Range("G19").Select()

Write-Host "Added rows 11-18 to sheet" $ws.Name
